# Update the EPEX Spot price workbook:
#  - "Prix Spot": insert a new date column (10-nov) before the 01-oct. block,
#    shifting the October/older columns one to the right, and stamp the new
#    column with "-" placeholders for every data row.
#  - "Gaz" and "CO2": append two new daily rows (2025-11-08, 2025-11-09)
#    carrying forward the last known price.

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert new "10-nov" column before DM ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Range("DM1").EntireColumn.Insert()
$wsSpot.Range("DM1").Value = "10-nov"

for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 117).Value = "-"
}

# --- Sheet "Gaz": append 2025-11-08 / 2025-11-09, carrying the last price ---
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("Z1").Value = "'2025-11-08"
$wsGaz.Range("Z1").Copy()
$wsGaz.Range("A146").PasteSpecial(-4163)
$wsGaz.Range("B146").Value = 29.755

$wsGaz.Range("Z1").Value = "'2025-11-09"
$wsGaz.Range("Z1").Copy()
$wsGaz.Range("A147").PasteSpecial(-4163)
$wsGaz.Range("B147").Value = 29.755

$wsGaz.Range("Z1").Clear()

# --- Sheet "CO2": append 2025-11-08 / 2025-11-09, carrying the last price ---
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("Z1").Value = "'2025-11-08"
$wsCo2.Range("Z1").Copy()
$wsCo2.Range("A146").PasteSpecial(-4163)
$wsCo2.Range("B146").Value = 79.36

$wsCo2.Range("Z1").Value = "'2025-11-09"
$wsCo2.Range("Z1").Copy()
$wsCo2.Range("A147").PasteSpecial(-4163)
$wsCo2.Range("B147").Value = 79.36

$wsCo2.Range("Z1").Clear()
